$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.936773538589478
$ws.Range("B1").Value = 2.014461278915405
$ws.Range("C1").Value = 1.915081739425659
$ws.Range("D1").Value = 1.050179004669189
$ws.Range("E1").Value = 0.7183305025100708
